# "Generate Report for Handoff" - update localization-status.xlsx after a new
# handoff package (b.*.xlf) was generated for e2e\b.md.
#
# b.md moves from "Handed back: in sync with en-US" to "Ready for handoff" in
# both target languages, a fresh handoff xliff is recorded (with a new
# "not the latest" warning because the source changed since the last
# handback), and the Overview rollup / column width follow along.

$wb = $excel.ActiveWorkbook

$statusReady  = "Ready for handoff"
$newHandoffDt = "2016-10-20 08:47:52"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7ef20be3e2f654b3101571fcbd2dc0e2ab400646/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f9b4981073477b48968a6e4cbf3b83dfa5f2dc/e2e/b.md."

# --- Overview sheet: b.md row (row 3) -------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusReady      # zh-cn status
$overview.Range("F3").Value = $statusReady      # de-de status
$overview.Range("G3").Value = $newHandoffDt     # Latest HO Xliff Generate Date

# --- zh-cn sheet: b.md row (row 3) -----------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusReady
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-10-20 08:47:40"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664

# --- de-de sheet: b.md row (row 3) -----------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusReady
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $newHandoffDt
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664
